$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UtilTimeMul")

# Update time multiplier samples in column B (rows 3-25) to follow the
# general load pattern instead of the flat block values.
$ws.Range("B3").Value = 0.78
$ws.Range("B4").Value = 0.72
$ws.Range("B5").Value = 0.71
$ws.Range("B6").Value = 0.76
$ws.Range("B7").Value = 0.97
$ws.Range("B8").Value = 1.1000000000000001
$ws.Range("B9").Value = 1.2
$ws.Range("B10").Value = 0.95
$ws.Range("B11").Value = 0.91
$ws.Range("B12").Value = 0.86
$ws.Range("B13").Value = 0.82
$ws.Range("B14").Value = 0.8
$ws.Range("B15").Value = 0.94
$ws.Range("B16").Value = 0.95
$ws.Range("B17").Value = 1.03
$ws.Range("B18").Value = 1.19
$ws.Range("B19").Value = 1.25
$ws.Range("B20").Value = 1.23
$ws.Range("B21").Value = 1.18
$ws.Range("B22").Value = 1.03
$ws.Range("B23").Value = 0.99
$ws.Range("B24").Value = 0.92
$ws.Range("B25").Value = 0.84

# Update the sheet's selected cell to match the new active cell/selection.
$ws.Activate()
$ws.Range("E17").Select()
